$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Insert two new rows right after the header row (row 1), pushing the
# existing 19 days of data (and the totals row) down by two rows.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Copy the date-format styling from the (now shifted) first data row so the
# two new date cells look like the rest of column A.
$ws.Range("A4").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row for 21 May 2024 (serial 45433)
$ws.Range("A2").Value = 45433
$ws.Range("B2").Value = 230
$ws.Range("C2").Value = 25
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 202

# New row for 20 May 2024 (serial 45432)
$ws.Range("A3").Value = 45432
$ws.Range("B3").Value = 232
$ws.Range("C3").Value = 24
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 205

# Fix up the totals row (now on row 23) so it sums the full, expanded range.
$ws.Range("B23").Formula = "=SUM(B2:B22)"
$ws.Range("C23").Formula = "=SUM(C2:C22)"
$ws.Range("D23").Formula = "=SUM(D2:D22)"
$ws.Range("E23").Formula = "=SUM(E2:E22)"

# Restore the cursor/selection position as left by the author.
$ws.Range("H22").Select()
